$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.472738
$ws.Range("H2").Value = 1.418214
$ws.Range("I2").Value = 0.0327564895931267
$ws.Range("J2").Value = 0.03397138804734427
$ws.Range("M2").Value = 26.23985166666667
$ws.Range("N2").Value = 78.719555
$ws.Range("O2").Value = 0.09560625159817936
$ws.Range("P2").Value = 0.09631575414263932
$ws.Range("Q2").Value = 12.40457499719667
$ws.Range("R2").Value = 111.64117497477
$ws.Range("S2").Value = 0.003131725185513615
$ws.Range("T2").Value = 0.003271979859052207
$ws.Range("G3").Value = 0.472738
$ws.Range("H3").Value = 1.418214
$ws.Range("I3").Value = 0.0327564895931267
$ws.Range("J3").Value = 0.03397138804734427
$ws.Range("O3").Value = 0.1052038484825964
$ws.Range("P3").Value = 0.1059845756519775
$ws.Range("Q3").Value = 13.649829448192
$ws.Range("R3").Value = 122.848465033728
$ws.Range("S3").Value = 0.003446108767977048
$ws.Range("T3").Value = 0.003600443146506443
$ws.Range("G4").Value = 0.472738
$ws.Range("H4").Value = 1.418214
$ws.Range("I4").Value = 0.0327564895931267
$ws.Range("J4").Value = 0.03397138804734427
$ws.Range("M4").Value = 108.455335
$ws.Range("N4").Value = 325.366005
$ws.Range("O4").Value = 0.3951626014085634
$ws.Range("P4").Value = 0.3980951384183098
$ws.Range("Q4").Value = 51.27095815723
$ws.Range("R4").Value = 461.43862341507
$ws.Range("S4").Value = 0.01294413964063248
$ws.Range("T4").Value = 0.01352384442696963
$ws.Range("G5").Value = 0.472738
$ws.Range("H5").Value = 1.418214
$ws.Range("I5").Value = 0.0327564895931267
$ws.Range("J5").Value = 0.03397138804734427
$ws.Range("M5").Value = 6.0653095
$ws.Range("N5").Value = 12.130619
$ws.Range("O5").Value = 0.02209926768810472
$ws.Range("P5").Value = 0.01484217888683478
$ws.Range("Q5").Value = 2.867302282411
$ws.Range("R5").Value = 17.203813694466
$ws.Range("S5").Value = 0.0007238944320411234
$ws.Range("T5").Value = 0.0005042094184327645
$ws.Range("G6").Value = 0.472738
$ws.Range("H6").Value = 1.418214
$ws.Range("I6").Value = 0.0327564895931267
$ws.Range("J6").Value = 0.03397138804734427
$ws.Range("M6").Value = 104.823008
$ws.Range("N6").Value = 314.469024
$ws.Range("O6").Value = 0.3819280308225561
$ws.Range("P6").Value = 0.3847623529002386
$ws.Range("Q6").Value = 49.55381915590401
$ws.Range("R6").Value = 445.984372403136
$ws.Range("S6").Value = 0.01251062156696243
$ws.Range("T6").Value = 0.01307091119638322
$ws.Range("I7").Value = 0.822180234441485
$ws.Range("J7").Value = 0.8526739017519405
$ws.Range("M7").Value = 26.23985166666667
$ws.Range("N7").Value = 78.719555
$ws.Range("O7").Value = 0.09560625159817936
$ws.Range("P7").Value = 0.09631575414263932
$ws.Range("Q7").Value = 311.3519338006889
$ws.Range("R7").Value = 2802.1674042062
$ws.Range("S7").Value = 0.0786055703530627
$ws.Range("T7").Value = 0.08212592988498489
$ws.Range("I8").Value = 0.822180234441485
$ws.Range("J8").Value = 0.8526739017519405
$ws.Range("O8").Value = 0.1052038484825964
$ws.Range("P8").Value = 0.1059845756519775
$ws.Range("S8").Value = 0.08649652480956758
$ws.Range("T8").Value = 0.09037028164669536
$ws.Range("I9").Value = 0.822180234441485
$ws.Range("J9").Value = 0.8526739017519405
$ws.Range("M9").Value = 108.455335
$ws.Range("N9").Value = 325.366005
$ws.Range("O9").Value = 0.3951626014085634
$ws.Range("P9").Value = 0.3980951384183098
$ws.Range("Q9").Value = 1286.889069047133
$ws.Range("R9").Value = 11582.0016214242
$ws.Range("S9").Value = 0.3248948802685997
$ws.Range("T9").Value = 0.339445334943619
$ws.Range("I10").Value = 0.822180234441485
$ws.Range("J10").Value = 0.8526739017519405
$ws.Range("M10").Value = 6.0653095
$ws.Range("N10").Value = 12.130619
$ws.Range("O10").Value = 0.02209926768810472
$ws.Range("P10").Value = 0.01484217888683478
$ws.Range("Q10").Value = 71.96861727399333
$ws.Range("R10").Value = 431.81170364396
$ws.Range("S10").Value = 0.01816958108879107
$ws.Range("T10").Value = 0.01265553858193768
$ws.Range("I11").Value = 0.822180234441485
$ws.Range("J11").Value = 0.8526739017519405
$ws.Range("M11").Value = 104.823008
$ws.Range("N11").Value = 314.469024
$ws.Range("O11").Value = 0.3819280308225561
$ws.Range("P11").Value = 0.3847623529002386
$ws.Range("Q11").Value = 1243.789281364907
$ws.Range("R11").Value = 11194.10353228416
$ws.Range("S11").Value = 0.3140136779214639
$ws.Range("T11").Value = 0.3280768166947035
$ws.Range("G12").Value = 0.37892
$ws.Range("H12").Value = 1.13676
$ws.Range("I12").Value = 0.02625574638939025
$ws.Range("J12").Value = 0.02722954016579943
$ws.Range("M12").Value = 26.23985166666667
$ws.Range("N12").Value = 78.719555
$ws.Range("O12").Value = 0.09560625159817936
$ws.Range("P12").Value = 0.09631575414263932
$ws.Range("Q12").Value = 9.942804593533333
$ws.Range("R12").Value = 89.48524134180001
$ws.Range("S12").Value = 0.002510213495202033
$ws.Range("T12").Value = 0.00262263369602626
$ws.Range("G13").Value = 0.37892
$ws.Range("H13").Value = 1.13676
$ws.Range("I13").Value = 0.02625574638939025
$ws.Range("J13").Value = 0.02722954016579943
$ws.Range("O13").Value = 0.1052038484825964
$ws.Range("P13").Value = 0.1059845756519775
$ws.Range("Q13").Value = 10.94093001728
$ws.Range("R13").Value = 98.46837015551999
$ws.Range("S13").Value = 0.002762205564946889
$ws.Range("T13").Value = 0.002885911259670729
$ws.Range("G14").Value = 0.37892
$ws.Range("H14").Value = 1.13676
$ws.Range("I14").Value = 0.02625574638939025
$ws.Range("J14").Value = 0.02722954016579943
$ws.Range("M14").Value = 108.455335
$ws.Range("N14").Value = 325.366005
$ws.Range("O14").Value = 0.3951626014085634
$ws.Range("P14").Value = 0.3980951384183098
$ws.Range("Q14").Value = 41.0958955382
$ws.Range("R14").Value = 369.8630598438
$ws.Range("S14").Value = 0.01037528904515495
$ws.Range("T14").Value = 0.01083994756137085
$ws.Range("G15").Value = 0.37892
$ws.Range("H15").Value = 1.13676
$ws.Range("I15").Value = 0.02625574638939025
$ws.Range("J15").Value = 0.02722954016579943
$ws.Range("M15").Value = 6.0653095
$ws.Range("N15").Value = 12.130619
$ws.Range("O15").Value = 0.02209926768810472
$ws.Range("P15").Value = 0.01484217888683478
$ws.Range("Q15").Value = 2.29826707574
$ws.Range("R15").Value = 13.78960245444
$ws.Range("S15").Value = 0.000580232767810124
$ws.Range("T15").Value = 0.0004041457061470479
$ws.Range("G16").Value = 0.37892
$ws.Range("H16").Value = 1.13676
$ws.Range("I16").Value = 0.02625574638939025
$ws.Range("J16").Value = 0.02722954016579943
$ws.Range("M16").Value = 104.823008
$ws.Range("N16").Value = 314.469024
$ws.Range("O16").Value = 0.3819280308225561
$ws.Range("P16").Value = 0.3847623529002386
$ws.Range("Q16").Value = 39.71953419136
$ws.Range("R16").Value = 357.47580772224
$ws.Range("S16").Value = 0.01002780551627626
$ws.Range("T16").Value = 0.01047690194258454
$ws.Range("G17").Value = 1.548357
$ws.Range("H17").Value = 3.096714
$ws.Range("I17").Value = 0.1072872076222874
$ws.Range("J17").Value = 0.0741775733180209
$ws.Range("M17").Value = 26.23985166666667
$ws.Range("N17").Value = 78.719555
$ws.Range("O17").Value = 0.09560625159817936
$ws.Range("P17").Value = 0.09631575414263932
$ws.Range("Q17").Value = 40.628658007045
$ws.Range("R17").Value = 243.77194804227
$ws.Range("S17").Value = 0.01025732776520251
$ws.Range("T17").Value = 0.007144468914596103
$ws.Range("G18").Value = 1.548357
$ws.Range("H18").Value = 3.096714
$ws.Range("I18").Value = 0.1072872076222874
$ws.Range("J18").Value = 0.0741775733180209
$ws.Range("O18").Value = 0.1052038484825964
$ws.Range("P18").Value = 0.1059845756519775
$ws.Range("Q18").Value = 44.707235244288
$ws.Range("R18").Value = 268.243411465728
$ws.Range("S18").Value = 0.01128702713481598
$ws.Range("T18").Value = 0.007861678631003894
$ws.Range("G19").Value = 1.548357
$ws.Range("H19").Value = 3.096714
$ws.Range("I19").Value = 0.1072872076222874
$ws.Range("J19").Value = 0.0741775733180209
$ws.Range("M19").Value = 108.455335
$ws.Range("N19").Value = 325.366005
$ws.Range("O19").Value = 0.3951626014085634
$ws.Range("P19").Value = 0.3980951384183098
$ws.Range("Q19").Value = 167.927577134595
$ws.Range("R19").Value = 1007.56546280757
$ws.Range("S19").Value = 0.04239589206188372
$ws.Range("T19").Value = 0.02952973131757185
$ws.Range("G20").Value = 1.548357
$ws.Range("H20").Value = 3.096714
$ws.Range("I20").Value = 0.1072872076222874
$ws.Range("J20").Value = 0.0741775733180209
$ws.Range("M20").Value = 6.0653095
$ws.Range("N20").Value = 12.130619
$ws.Range("O20").Value = 0.02209926768810472
$ws.Range("P20").Value = 0.01484217888683478
$ws.Range("Q20").Value = 9.391264421491501
$ws.Range("R20").Value = 37.565057685966
$ws.Range("S20").Value = 0.002370968720754197
$ws.Range("T20").Value = 0.001100956812577369
$ws.Range("G21").Value = 1.548357
$ws.Range("H21").Value = 3.096714
$ws.Range("I21").Value = 0.1072872076222874
$ws.Range("J21").Value = 0.0741775733180209
$ws.Range("M21").Value = 104.823008
$ws.Range("N21").Value = 314.469024
$ws.Range("O21").Value = 0.3819280308225561
$ws.Range("P21").Value = 0.3847623529002386
$ws.Range("Q21").Value = 162.303438197856
$ws.Range("R21").Value = 973.8206291871361
$ws.Range("S21").Value = 0.04097599193963094
$ws.Range("T21").Value = 0.02854073764227168
$ws.Range("G22").Value = 0.16626
$ws.Range("H22").Value = 0.49878
$ws.Range("I22").Value = 0.01152032195371061
$ws.Range("J22").Value = 0.01194759671689489
$ws.Range("M22").Value = 26.23985166666667
$ws.Range("N22").Value = 78.719555
$ws.Range("O22").Value = 0.09560625159817936
$ws.Range("P22").Value = 0.09631575414263932
$ws.Range("Q22").Value = 4.3626377381
$ws.Range("R22").Value = 39.2637396429
$ws.Range("S22").Value = 0.001101414799198485
$ws.Range("T22").Value = 0.001150741787979853
$ws.Range("G23").Value = 0.16626
$ws.Range("H23").Value = 0.49878
$ws.Range("I23").Value = 0.01152032195371061
$ws.Range("J23").Value = 0.01194759671689489
$ws.Range("O23").Value = 0.1052038484825964
$ws.Range("P23").Value = 0.1059845756519775
$ws.Range("Q23").Value = 4.800588579839999
$ws.Range("R23").Value = 43.20529721856
$ws.Range("S23").Value = 0.0012119822052889
$ws.Range("T23").Value = 0.001266260968101065
$ws.Range("G24").Value = 0.16626
$ws.Range("H24").Value = 0.49878
$ws.Range("I24").Value = 0.01152032195371061
$ws.Range("J24").Value = 0.01194759671689489
$ws.Range("M24").Value = 108.455335
$ws.Range("N24").Value = 325.366005
$ws.Range("O24").Value = 0.3951626014085634
$ws.Range("P24").Value = 0.3980951384183098
$ws.Range("Q24").Value = 18.0317839971
$ws.Range("R24").Value = 162.2860559739
$ws.Range("S24").Value = 0.004552400392292467
$ws.Range("T24").Value = 0.004756280168778416
$ws.Range("G25").Value = 0.16626
$ws.Range("H25").Value = 0.49878
$ws.Range("I25").Value = 0.01152032195371061
$ws.Range("J25").Value = 0.01194759671689489
$ws.Range("M25").Value = 6.0653095
$ws.Range("N25").Value = 12.130619
$ws.Range("O25").Value = 0.02209926768810472
$ws.Range("P25").Value = 0.01484217888683478
$ws.Range("Q25").Value = 1.00841835747
$ws.Range("R25").Value = 6.05051014482
$ws.Range("S25").Value = 0.0002545906787082002
$ws.Range("T25").Value = 0.0001773283677399139
$ws.Range("G26").Value = 0.16626
$ws.Range("H26").Value = 0.49878
$ws.Range("I26").Value = 0.01152032195371061
$ws.Range("J26").Value = 0.01194759671689489
$ws.Range("M26").Value = 104.823008
$ws.Range("N26").Value = 314.469024
$ws.Range("O26").Value = 0.3819280308225561
$ws.Range("P26").Value = 0.3847623529002386
$ws.Range("Q26").Value = 17.42787331008
$ws.Range("R26").Value = 156.85085979072
$ws.Range("S26").Value = 0.004399933878222555
